$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header row: "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404"
#    (the "diff" header in column K is unchanged)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# ---------------------------------------------------------------------------
# 2. Turn the data range into an Excel Table ("Table1"), spanning the full
#    used range A1:U72. The header row already carries a bold/filled/bordered
#    style (style index 1) baked into the sheet; stash a copy of that format
#    on a scratch cell, clear the header formatting so ListObjects.Add()
#    doesn't promote it into a table "headerRowDxfId" override, add the
#    table, then paste the original formatting back onto the header row.
# ---------------------------------------------------------------------------
$ws.Range("A1").Copy() | Out-Null
$ws.Range("W1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = $false

$headerRange = $ws.Range("A1:U1")
$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U72"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""
$tbl.ShowTableStyleRowStripes = $true
$tbl.ShowTableStyleColumnStripes = $false
$tbl.ShowTableStyleFirstColumn = $false
$tbl.ShowTableStyleLastColumn = $false

$ws.Range("W1").Copy() | Out-Null
$headerRange.PasteSpecial(-4122) | Out-Null      # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("W1").Clear()

# ---------------------------------------------------------------------------
# 3. Freeze the header row (row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
